# Auto-generated Excel COM-interop script applying the 2026-01-30 16:18:44 scrape refresh
# for horarios-141 workbook (3 sheets: LP1912, LP1912-215, 6203-6173).
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 16:18:44'
$ws.Cells.Item(3, 1).Value = 'Total filas: 273'
$ws.Cells.Item(38, 1).Value = '07:23:38'
$ws.Cells.Item(38, 2).Value = '07:36'
$ws.Cells.Item(38, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(38, 4).Value = 13
$ws.Cells.Item(38, 5).Value = 'LP1912'
$ws.Cells.Item(39, 1).Value = '07:23:38'
$ws.Cells.Item(39, 2).Value = '07:36'
$ws.Cells.Item(39, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(39, 4).Value = 13
$ws.Cells.Item(39, 5).Value = 'LP1912'
$ws.Cells.Item(60, 1).Value = '07:23:38'
$ws.Cells.Item(60, 2).Value = '08:44'
$ws.Cells.Item(60, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(60, 4).Value = 81
$ws.Cells.Item(60, 5).Value = 'LP1912'
$ws.Cells.Item(61, 1).Value = '08:31:16'
$ws.Cells.Item(61, 2).Value = '08:44'
$ws.Cells.Item(61, 3).Value = '14_ABASTO'
$ws.Cells.Item(61, 4).Value = 13
$ws.Cells.Item(61, 5).Value = 'LP1912'
$ws.Cells.Item(88, 1).Value = '09:32:47'
$ws.Cells.Item(88, 2).Value = '09:42'
$ws.Cells.Item(88, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(88, 4).Value = 10
$ws.Cells.Item(88, 5).Value = 'LP1912'
$ws.Cells.Item(89, 1).Value = '08:31:16'
$ws.Cells.Item(89, 2).Value = '09:42'
$ws.Cells.Item(89, 3).Value = '215C_EL PATO'
$ws.Cells.Item(89, 4).Value = 71
$ws.Cells.Item(89, 5).Value = 'LP1912'
$ws.Cells.Item(112, 1).Value = '10:39:14'
$ws.Cells.Item(112, 2).Value = '10:57'
$ws.Cells.Item(112, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(112, 4).Value = 18
$ws.Cells.Item(112, 5).Value = 'LP1912'
$ws.Cells.Item(113, 1).Value = '10:39:14'
$ws.Cells.Item(113, 2).Value = '10:57'
$ws.Cells.Item(113, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(113, 4).Value = 18
$ws.Cells.Item(113, 5).Value = 'LP1912'
$ws.Cells.Item(152, 1).Value = '11:57:34'
$ws.Cells.Item(152, 2).Value = '12:21'
$ws.Cells.Item(152, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(152, 4).Value = 24
$ws.Cells.Item(152, 5).Value = 'LP1912'
$ws.Cells.Item(153, 1).Value = '10:39:14'
$ws.Cells.Item(153, 2).Value = '12:21'
$ws.Cells.Item(153, 3).Value = '215A_EL PATO'
$ws.Cells.Item(153, 4).Value = 102
$ws.Cells.Item(153, 5).Value = 'LP1912'
$ws.Cells.Item(159, 1).Value = '12:30:50'
$ws.Cells.Item(159, 2).Value = '12:36'
$ws.Cells.Item(159, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(159, 4).Value = 6
$ws.Cells.Item(159, 5).Value = 'LP1912'
$ws.Cells.Item(160, 1).Value = '12:30:50'
$ws.Cells.Item(160, 2).Value = '12:36'
$ws.Cells.Item(160, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(160, 4).Value = 6
$ws.Cells.Item(160, 5).Value = 'LP1912'
$ws.Cells.Item(177, 1).Value = '11:57:34'
$ws.Cells.Item(177, 2).Value = '13:21'
$ws.Cells.Item(177, 3).Value = '10_OLMOS'
$ws.Cells.Item(177, 4).Value = 84
$ws.Cells.Item(177, 5).Value = 'LP1912'
$ws.Cells.Item(178, 1).Value = '13:02:37'
$ws.Cells.Item(178, 2).Value = '13:21'
$ws.Cells.Item(178, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(178, 4).Value = 19
$ws.Cells.Item(178, 5).Value = 'LP1912'
$ws.Cells.Item(185, 1).Value = '13:02:37'
$ws.Cells.Item(185, 2).Value = '13:46'
$ws.Cells.Item(185, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(185, 4).Value = 44
$ws.Cells.Item(185, 5).Value = 'LP1912'
$ws.Cells.Item(186, 1).Value = '13:02:37'
$ws.Cells.Item(186, 2).Value = '13:46'
$ws.Cells.Item(186, 3).Value = '17_ROMERO'
$ws.Cells.Item(186, 4).Value = 44
$ws.Cells.Item(186, 5).Value = 'LP1912'
$ws.Cells.Item(190, 1).Value = '13:02:37'
$ws.Cells.Item(190, 2).Value = '13:56'
$ws.Cells.Item(190, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(190, 4).Value = 54
$ws.Cells.Item(190, 5).Value = 'LP1912'
$ws.Cells.Item(191, 1).Value = '13:02:37'
$ws.Cells.Item(191, 2).Value = '13:56'
$ws.Cells.Item(191, 3).Value = '225_GOMEZ'
$ws.Cells.Item(191, 4).Value = 54
$ws.Cells.Item(191, 5).Value = 'LP1912'
$ws.Cells.Item(207, 1).Value = '14:26:27'
$ws.Cells.Item(207, 2).Value = '14:44'
$ws.Cells.Item(207, 3).Value = '15_ABASTO'
$ws.Cells.Item(207, 4).Value = 18
$ws.Cells.Item(207, 5).Value = 'LP1912'
$ws.Cells.Item(208, 1).Value = '14:26:27'
$ws.Cells.Item(208, 2).Value = '14:44'
$ws.Cells.Item(208, 3).Value = '14_ABASTO'
$ws.Cells.Item(208, 4).Value = 18
$ws.Cells.Item(208, 5).Value = 'LP1912'
$ws.Cells.Item(228, 1).Value = '14:57:41'
$ws.Cells.Item(228, 2).Value = '15:37'
$ws.Cells.Item(228, 3).Value = '10_OLMOS'
$ws.Cells.Item(228, 4).Value = 40
$ws.Cells.Item(228, 5).Value = 'LP1912'
$ws.Cells.Item(229, 1).Value = '15:37:31'
$ws.Cells.Item(229, 2).Value = '15:37'
$ws.Cells.Item(229, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(229, 4).Value = 0
$ws.Cells.Item(229, 5).Value = 'LP1912'
$ws.Cells.Item(230, 1).Value = '14:26:27'
$ws.Cells.Item(230, 2).Value = '15:38'
$ws.Cells.Item(230, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(230, 4).Value = 72
$ws.Cells.Item(230, 5).Value = 'LP1912'
$ws.Cells.Item(231, 1).Value = '15:37:31'
$ws.Cells.Item(231, 2).Value = '15:38'
$ws.Cells.Item(231, 3).Value = '215A_EL PATO'
$ws.Cells.Item(231, 4).Value = 1
$ws.Cells.Item(231, 5).Value = 'LP1912'
$ws.Cells.Item(240, 1).Value = '15:37:31'
$ws.Cells.Item(240, 2).Value = '15:54'
$ws.Cells.Item(240, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(240, 4).Value = 17
$ws.Cells.Item(240, 5).Value = 'LP1912'
$ws.Cells.Item(241, 1).Value = '13:55:05'
$ws.Cells.Item(241, 2).Value = '15:54'
$ws.Cells.Item(241, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(241, 4).Value = 119
$ws.Cells.Item(241, 5).Value = 'LP1912'
$ws.Cells.Item(248, 1).Value = '16:18:44'
$ws.Cells.Item(248, 2).Value = '16:19'
$ws.Cells.Item(248, 3).Value = '215C_EL PATO'
$ws.Cells.Item(248, 4).Value = 1
$ws.Cells.Item(248, 5).Value = 'LP1912'
$ws.Cells.Item(250, 1).Value = '16:18:44'
$ws.Cells.Item(250, 2).Value = '16:21'
$ws.Cells.Item(250, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(250, 4).Value = 3
$ws.Cells.Item(250, 5).Value = 'LP1912'
$ws.Cells.Item(251, 1).Value = '16:18:44'
$ws.Cells.Item(251, 2).Value = '16:26'
$ws.Cells.Item(251, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(251, 4).Value = 8
$ws.Cells.Item(251, 5).Value = 'LP1912'
$ws.Cells.Item(252, 1).Value = '16:18:44'
$ws.Cells.Item(252, 2).Value = '16:29'
$ws.Cells.Item(252, 3).Value = '10_OLMOS'
$ws.Cells.Item(252, 4).Value = 11
$ws.Cells.Item(252, 5).Value = 'LP1912'
$ws.Cells.Item(253, 1).Value = '16:18:44'
$ws.Cells.Item(253, 2).Value = '16:30'
$ws.Cells.Item(253, 3).Value = '15_ABASTO'
$ws.Cells.Item(253, 4).Value = 12
$ws.Cells.Item(253, 5).Value = 'LP1912'
$ws.Cells.Item(254, 1).Value = '16:18:44'
$ws.Cells.Item(254, 2).Value = '16:34'
$ws.Cells.Item(254, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(254, 4).Value = 16
$ws.Cells.Item(254, 5).Value = 'LP1912'
$ws.Cells.Item(255, 1).Value = '16:18:44'
$ws.Cells.Item(255, 2).Value = '16:34'
$ws.Cells.Item(255, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(255, 4).Value = 16
$ws.Cells.Item(255, 5).Value = 'LP1912'
$ws.Cells.Item(256, 1).Value = '15:37:31'
$ws.Cells.Item(256, 2).Value = '16:36'
$ws.Cells.Item(256, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(256, 4).Value = 59
$ws.Cells.Item(256, 5).Value = 'LP1912'
$ws.Cells.Item(257, 1).Value = '15:37:31'
$ws.Cells.Item(257, 2).Value = '16:38'
$ws.Cells.Item(257, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(257, 4).Value = 61
$ws.Cells.Item(257, 5).Value = 'LP1912'
$ws.Cells.Item(258, 1).Value = '15:37:31'
$ws.Cells.Item(258, 2).Value = '16:40'
$ws.Cells.Item(258, 3).Value = '17_ROMERO'
$ws.Cells.Item(258, 4).Value = 63
$ws.Cells.Item(258, 5).Value = 'LP1912'
$ws.Cells.Item(259, 1).Value = '16:18:44'
$ws.Cells.Item(259, 2).Value = '16:42'
$ws.Cells.Item(259, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(259, 4).Value = 24
$ws.Cells.Item(259, 5).Value = 'LP1912'
$ws.Cells.Item(260, 1).Value = '16:18:44'
$ws.Cells.Item(260, 2).Value = '16:43'
$ws.Cells.Item(260, 3).Value = '225_GOMEZ'
$ws.Cells.Item(260, 4).Value = 25
$ws.Cells.Item(260, 5).Value = 'LP1912'
$ws.Cells.Item(261, 1).Value = '16:18:44'
$ws.Cells.Item(261, 2).Value = '16:48'
$ws.Cells.Item(261, 3).Value = '15_ABASTO'
$ws.Cells.Item(261, 4).Value = 30
$ws.Cells.Item(261, 5).Value = 'LP1912'
$ws.Cells.Item(262, 1).Value = '16:18:44'
$ws.Cells.Item(262, 2).Value = '16:50'
$ws.Cells.Item(262, 3).Value = '14_ABASTO'
$ws.Cells.Item(262, 4).Value = 32
$ws.Cells.Item(262, 5).Value = 'LP1912'
$ws.Cells.Item(263, 1).Value = '16:18:44'
$ws.Cells.Item(263, 2).Value = '16:56'
$ws.Cells.Item(263, 3).Value = '17_179 Y 38'
$ws.Cells.Item(263, 4).Value = 38
$ws.Cells.Item(263, 5).Value = 'LP1912'
$ws.Cells.Item(264, 1).Value = '16:18:44'
$ws.Cells.Item(264, 2).Value = '16:57'
$ws.Cells.Item(264, 3).Value = '10_OLMOS'
$ws.Cells.Item(264, 4).Value = 39
$ws.Cells.Item(264, 5).Value = 'LP1912'
$ws.Cells.Item(265, 1).Value = '16:18:44'
$ws.Cells.Item(265, 2).Value = '17:04'
$ws.Cells.Item(265, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(265, 4).Value = 46
$ws.Cells.Item(265, 5).Value = 'LP1912'
$ws.Cells.Item(266, 1).Value = '16:18:44'
$ws.Cells.Item(266, 2).Value = '17:04'
$ws.Cells.Item(266, 3).Value = '215A_EL PATO'
$ws.Cells.Item(266, 4).Value = 46
$ws.Cells.Item(266, 5).Value = 'LP1912'
$ws.Cells.Item(267, 1).Value = '16:18:44'
$ws.Cells.Item(267, 2).Value = '17:04'
$ws.Cells.Item(267, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(267, 4).Value = 46
$ws.Cells.Item(267, 5).Value = 'LP1912'
$ws.Cells.Item(268, 1).Value = '16:18:44'
$ws.Cells.Item(268, 2).Value = '17:16'
$ws.Cells.Item(268, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(268, 4).Value = 58
$ws.Cells.Item(268, 5).Value = 'LP1912'
$ws.Cells.Item(269, 1).Value = '16:18:44'
$ws.Cells.Item(269, 2).Value = '17:21'
$ws.Cells.Item(269, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(269, 4).Value = 63
$ws.Cells.Item(269, 5).Value = 'LP1912'
$ws.Cells.Item(270, 1).Value = '16:18:44'
$ws.Cells.Item(270, 2).Value = '17:24'
$ws.Cells.Item(270, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(270, 4).Value = 66
$ws.Cells.Item(270, 5).Value = 'LP1912'
$ws.Cells.Item(271, 1).Value = '16:18:44'
$ws.Cells.Item(271, 2).Value = '17:28'
$ws.Cells.Item(271, 3).Value = '14_ABASTO'
$ws.Cells.Item(271, 4).Value = 70
$ws.Cells.Item(271, 5).Value = 'LP1912'
$ws.Cells.Item(272, 1).Value = '16:18:44'
$ws.Cells.Item(272, 2).Value = '17:32'
$ws.Cells.Item(272, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(272, 4).Value = 74
$ws.Cells.Item(272, 5).Value = 'LP1912'
$ws.Cells.Item(273, 1).Value = '15:37:31'
$ws.Cells.Item(273, 2).Value = '17:36'
$ws.Cells.Item(273, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(273, 4).Value = 119
$ws.Cells.Item(273, 5).Value = 'LP1912'
$ws.Cells.Item(274, 1).Value = '16:18:44'
$ws.Cells.Item(274, 2).Value = '17:38'
$ws.Cells.Item(274, 3).Value = '17_ROMERO'
$ws.Cells.Item(274, 4).Value = 80
$ws.Cells.Item(274, 5).Value = 'LP1912'
$ws.Cells.Item(275, 1).Value = '16:18:44'
$ws.Cells.Item(275, 2).Value = '17:40'
$ws.Cells.Item(275, 3).Value = '215B_EL PATO'
$ws.Cells.Item(275, 4).Value = 82
$ws.Cells.Item(275, 5).Value = 'LP1912'
$ws.Cells.Item(276, 1).Value = '16:18:44'
$ws.Cells.Item(276, 2).Value = '17:50'
$ws.Cells.Item(276, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(276, 4).Value = 92
$ws.Cells.Item(276, 5).Value = 'LP1912'
$ws.Cells.Item(277, 1).Value = '16:18:44'
$ws.Cells.Item(277, 2).Value = '17:52'
$ws.Cells.Item(277, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(277, 4).Value = 94
$ws.Cells.Item(277, 5).Value = 'LP1912'
$ws.Cells.Item(278, 1).Value = '16:18:44'
$ws.Cells.Item(278, 2).Value = '18:04'
$ws.Cells.Item(278, 3).Value = '17_ROMERO'
$ws.Cells.Item(278, 4).Value = 106
$ws.Cells.Item(278, 5).Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 16:18:44'
$ws.Cells.Item(3, 1).Value = 'Total filas: 28'
$ws.Cells.Item(30, 1).Value = '16:18:44'
$ws.Cells.Item(30, 2).Value = '16:19'
$ws.Cells.Item(30, 3).Value = '215C_EL PATO'
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(30, 5).Value = 'LP1912'
$ws.Cells.Item(32, 1).Value = '16:18:44'
$ws.Cells.Item(32, 2).Value = '17:04'
$ws.Cells.Item(32, 3).Value = '215A_EL PATO'
$ws.Cells.Item(32, 4).Value = 46
$ws.Cells.Item(32, 5).Value = 'LP1912'
$ws.Cells.Item(33, 1).Value = '16:18:44'
$ws.Cells.Item(33, 2).Value = '17:40'
$ws.Cells.Item(33, 3).Value = '215B_EL PATO'
$ws.Cells.Item(33, 4).Value = 82
$ws.Cells.Item(33, 5).Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 16:18:44'
$ws.Cells.Item(3, 1).Value = 'Total filas: 35'
$ws.Cells.Item(36, 1).Value = '16:18:44'
$ws.Cells.Item(36, 2).Value = '16:52'
$ws.Cells.Item(36, 3).Value = '215B_LP-P MOR-40 Y 115'
$ws.Cells.Item(36, 4).Value = 34
$ws.Cells.Item(36, 5).Value = 'L6173'
$ws.Cells.Item(39, 1).Value = '16:18:44'
$ws.Cells.Item(39, 2).Value = '17:20'
$ws.Cells.Item(39, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(39, 4).Value = 62
$ws.Cells.Item(39, 5).Value = 'L6173'
$ws.Cells.Item(40, 1).Value = '16:18:44'
$ws.Cells.Item(40, 2).Value = '18:03'
$ws.Cells.Item(40, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(40, 4).Value = 105
$ws.Cells.Item(40, 5).Value = 'L6203'

